$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new pokemon entry "charmander" in column A, row 3 (below existing "totodile" row)
$ws.Range("A3").Value = "charmander"

# Update the active selection to A4, matching the saved cursor position in the file
$ws.Range("A4").Select()
